# Apply updates to "Mental Health Ontology mapping to LSRs.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple label edits (single cell F column updates) ---
$ws.Range("F32").Value = "Mean drug dose delivered (milligrams)"
$ws.Range("F37").Value = "Number of participants reporting post-intervention constipation"
$ws.Range("F38").Value = "Number of participants reporting post-intervention dizziness"
$ws.Range("F39").Value = "Maximum drug dose delivered (milligrams)"
$ws.Range("F40").Value = "Mean drug dose delivered (milligrams)"
$ws.Range("F41").Value = "Minimum drug dose delivered (milligrams)"
$ws.Range("F43").Value = "Number of participants reporting post-intervention dry mouth"
$ws.Range("F47").Value = "Number of participants reporting post-intervention headache"
$ws.Range("F51").Value = "Number of participants reporting post-intervention insomnia"
$ws.Range("F55").Value = "Maximum drug dose delivered (milligrams)"
$ws.Range("F56").Value = "Maximum planned drug dose (milligrams)"
$ws.Range("F91").Value = "Minimum drug dose delivered (milligrams)"
$ws.Range("F92").Value = "Minimum planned drug dose (milligrams)"
$ws.Range("F93").Value = "Number of participants reporting post-intervention nausea"
$ws.Range("F101").Value = "Number of participants randomly allocated to an arm"
$ws.Range("F102").Value = "Maximum planned drug dose (milligrams)"
$ws.Range("F103").Value = "Minimum planned drug dose (milligrams)"
$ws.Range("F132").Value = "Number of participants reporting post-intervention vomiting"
$ws.Range("F137").Value = "Mean drug dose delivered (milligrams)"
$ws.Range("F138").Value = "Mean drug dose delivered (milligrams)"
$ws.Range("F331").Value = "Number of participants reporting post-intervention nausea"
$ws.Range("F332").Value = "Number of participants reporting post-intervention nausea"
$ws.Range("F333").Value = "Number of participants reporting post-intervention nausea"
$ws.Range("F615").Value = "Number of participants reporting post-intervention dizziness"
$ws.Range("F616").Value = "Number of participants reporting post-intervention dry mouth"
$ws.Range("F618").Value = "Number of participants reporting post-intervention headache"
$ws.Range("F619").Value = "Number of participants reporting post-intervention insomnia"
$ws.Range("F620").Value = "Number of participants reporting post-intervention nausea"
$ws.Range("F625").Value = "Number of participants reporting post-intervention vomiting"
$ws.Range("F627").Value = "Number of participants reporting post-intervention constipation"

# --- Row 35: full row re-mapping (participant drop-out) ---
$ws.Range("F35").Value = "Number of participants who dropped out of study"
$ws.Range("G35").Value = "GMHO:0000075"
$ws.Range("H35").Value = "number of participant drop-out from intervention"
$ws.Range("I35").Value = "Number of intervention participants who withdraw from or cannot complete an intervention."
$ws.Range("J35").Value = "number of intervention participants"
$ws.Range("L35").Value = "GMHO:0000152"

# --- Row 36: full row re-mapping (drop-out due to adverse events) ---
$ws.Range("F36").Value = "Number of participants who dropped out of study due to adverse events"
$ws.Range("G36").Value = "GMHO:0000072"
$ws.Range("H36").Value = "number of participant drop-out due to adverse events"
$ws.Range("I36").Value = "Number of participant drop-out from the intervention as a result the participants experiencing some adverse event."
$ws.Range("J36").Value = "number of participant drop-out from intervention"
$ws.Range("L36").Value = "GMHO:0000075"

# --- Row 575: full row re-mapping (number of intervention participants) ---
$ws.Range("F575").Value = "number of intervention participants"
$ws.Range("G575").Value = "GMHO:0000152"
$ws.Range("H575").Value = "number of intervention participants"
$ws.Range("I575").Value = "A data item that is about the number of participants in an intervention or part of an intervention at a timepoint."
$ws.Range("J575").Value = "data item"
$ws.Range("L575").Value = "Population"

# --- New rows 638-644 appended at the bottom of the sheet ---
$newRows = @(
    @{ E = "constipation"; F = "Number of participants reporting post-intervention constipation"; L = "GMHO:0000257,GMHO:0000204" },
    @{ E = "nausea";       F = "Number of participants reporting post-intervention nausea";       L = "GMHO:0000183,GMHO:0000204" },
    @{ E = "vomiting";     F = "Number of participants reporting post-intervention vomiting";     L = "GMHO:0000185,GMHO:0000204" },
    @{ E = "dizziness";    F = "Number of participants reporting post-intervention dizziness";    L = "GMHO:0000179,GMHO:0000204" },
    @{ E = "insomnia";     F = "Number of participants reporting post-intervention insomnia";     L = "GMHO:0000182,GMHO:0000204" },
    @{ E = "headache";     F = "Number of participants reporting post-intervention headache";     L = "GMHO:0000181,GMHO:0000204" },
    @{ E = "dry_mouth";    F = "Number of participants reporting post-intervention dry mouth";    L = "GMHO:0000180,GMHO:0000204" }
)

$rowNum = 638
foreach ($item in $newRows) {
    $ws.Range("A$rowNum").Value = "LSR1_anhedonia_H/data/human/df_amended_20240430.csv"
    # Column C holds the text "1" - force text storage so it is not coerced to a number.
    $ws.Range("C$rowNum").NumberFormat = "@"
    $ws.Range("C$rowNum").Value = "1"
    $ws.Range("D$rowNum").Value = "Outcomes"
    $ws.Range("E$rowNum").Value = $item.E
    $ws.Range("F$rowNum").Value = $item.F
    $ws.Range("G$rowNum").Value = "GMHO:0000214"
    $ws.Range("H$rowNum").Value = "number of participants with specific outcome"
    $ws.Range("I$rowNum").Value = "Number of intervention participants for whom an outcome has been identified within a study."
    $ws.Range("J$rowNum").Value = "number of intervention participants"
    $ws.Range("K$rowNum").Value = "Intervention outcomes and spillover effects"
    $ws.Range("L$rowNum").Value = $item.L
    $ws.Range("M$rowNum").Value = "No Combo"
    $rowNum++
}
